$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# New column header
$ws.Range("C1").Value = "Minimum Load"

# Minimum Load values per technology row (rows 2-13)
$ws.Range("C2").Value = 0      # Wind Offshore
$ws.Range("C3").Value = 0      # Wind Onshore
$ws.Range("C4").Value = 0      # Biomass
$ws.Range("C5").Value = 0.3    # Lignite
$ws.Range("C6").Value = 0.3    # Hard Coal
$ws.Range("C7").Value = 0      # Gas
$ws.Range("C8").Value = 0      # Hydrogen
$ws.Range("C9").Value = 0.3    # Nuclear
$ws.Range("C10").Value = 0     # PV
$ws.Range("C11").Value = 0     # Hydro
$ws.Range("C12").Value = 0     # Other RES
$ws.Range("C13").Value = 0     # Other Conventional

# Auto-size column C to fit its new "Minimum Load" content (mirrors the
# best-fit width the author ended up with after adding the column)
$ws.Columns.Item(3).ColumnWidth = 11.6

# Update the active selection to reflect the final cursor position
$ws.Range("C10").Select()
